# New weekly price report: insert a new record at the top of the data
# block (row 17, right after the most recent previously-known week),
# pushing all existing data rows down by one. This mirrors the original
# edit: a fresh week's quote was appended to the raw source and the
# report was regenerated, so every later row's values shift down by one
# position while a brand-new row is introduced with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 17:129 down to 18:130, leaving a blank row 17 to fill in.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row with this week's data.
$row = 17
$ws.Cells.Item($row, 1).Value  = 4
$ws.Cells.Item($row, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value  = "Los Lagos"
$ws.Cells.Item($row, 4).Value  = 44462
$ws.Cells.Item($row, 5).Value  = 10
$ws.Cells.Item($row, 6).Value  = 100112017
$ws.Cells.Item($row, 7).Value  = "Apio"
$ws.Cells.Item($row, 8).Value  = "Americana (o)"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 30
$ws.Cells.Item($row, 11).Value = 12000
$ws.Cells.Item($row, 12).Value = 12000
$ws.Cells.Item($row, 13).Value = 12000
$ws.Cells.Item($row, 14).Value = '$/docena de matas'
$ws.Cells.Item($row, 15).Value = "Región de Coquimbo"
$ws.Cells.Item($row, 16).Value = 2000
$ws.Cells.Item($row, 17).Value = 6
$ws.Cells.Item($row, 18).Value = "Hortaliza"
